# "Generate Report for Handoff"
# Refresh the handoff/handback timestamps shown in the localization-status
# report: the Overview sheet's "Latest Handoff Date" column, and the
# "Latest Handoff Datetime" column on each per-language sheet.

$wb = $excel.ActiveWorkbook

# Overview sheet: File Name / zh-cn / de-de / Latest Handoff Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D2:D5").Value = "2016-04-17 07:04:12"

# zh-cn sheet: column E is "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2:E5").Value = "2016-03-17 07:04:04"

# de-de sheet: column E is "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2:E5").Value = "2016-03-17 07:04:12"
